# Auto-generated cell updates applying the scheduled runner diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 112
$ws.Cells.Item(112, 8).Value = 8022456
$ws.Cells.Item(112, 10).Value = 8022456
$ws.Cells.Item(112, 12).Value = 24067368
$ws.Cells.Item(112, 14).Value = -24069584
# Row 135
$ws.Cells.Item(135, 8).Value = 1138.8235
$ws.Cells.Item(135, 9).Value = 1052.3334
$ws.Cells.Item(135, 10).Value = 1787.5
$ws.Cells.Item(135, 11).Value = 9471.000599999999
$ws.Cells.Item(135, 12).Value = 16087.5
$ws.Cells.Item(135, 13).Value = -6936.000599999999
$ws.Cells.Item(135, 14).Value = -21157.5
# Row 138
$ws.Cells.Item(138, 8).Value = 1800.0918
$ws.Cells.Item(138, 9).Value = 634
$ws.Cells.Item(138, 10).Value = 2419.5781
$ws.Cells.Item(138, 11).Value = 1902
$ws.Cells.Item(138, 12).Value = 7258.7343
$ws.Cells.Item(138, 13).Value = 3238
$ws.Cells.Item(138, 14).Value = -17538.7343

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Cells.Item(32, 8).Value = 2103.598
$ws.Cells.Item(32, 9).Value = 1404.3292
$ws.Cells.Item(32, 10).Value = 5926.2666
$ws.Cells.Item(32, 11).Value = 1404.3292
$ws.Cells.Item(32, 12).Value = 5926.2666
$ws.Cells.Item(32, 13).Value = -1117.3292
$ws.Cells.Item(32, 14).Value = -6500.2666
# Row 61
$ws.Cells.Item(61, 8).Value = 2936.9565
$ws.Cells.Item(61, 9).Value = 1863.6923
$ws.Cells.Item(61, 10).Value = 4332.2
$ws.Cells.Item(61, 11).Value = 1863.6923
$ws.Cells.Item(61, 12).Value = 4332.2
$ws.Cells.Item(61, 13).Value = -1651.6923
$ws.Cells.Item(61, 14).Value = -4756.2
# Row 74
$ws.Cells.Item(74, 8).Value = 1532.9286
$ws.Cells.Item(74, 9).Value = 1697.5294
$ws.Cells.Item(74, 11).Value = 1697.5294
$ws.Cells.Item(74, 13).Value = -823.5293999999999
# Row 77
$ws.Cells.Item(77, 8).Value = 1532.9286
$ws.Cells.Item(77, 9).Value = 1697.5294
$ws.Cells.Item(77, 11).Value = 8487.646999999999
$ws.Cells.Item(77, 13).Value = -4119.646999999999
# Row 108
$ws.Cells.Item(108, 8).Value = 30000
$ws.Cells.Item(108, 10).Value = 30000
$ws.Cells.Item(108, 12).Value = 30000
$ws.Cells.Item(108, 14).Value = -37680
# Row 136
$ws.Cells.Item(136, 8).Value = 2936.9565
$ws.Cells.Item(136, 9).Value = 1863.6923
$ws.Cells.Item(136, 10).Value = 4332.2
$ws.Cells.Item(136, 11).Value = 5591.0769
$ws.Cells.Item(136, 12).Value = 12996.6
$ws.Cells.Item(136, 13).Value = -3041.0769
$ws.Cells.Item(136, 14).Value = -18096.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Cells.Item(31, 8).Value = 1443.6538
$ws.Cells.Item(31, 9).Value = 1070.6
$ws.Cells.Item(31, 10).Value = 1747.3024
$ws.Cells.Item(31, 11).Value = 1070.6
$ws.Cells.Item(31, 12).Value = 1747.3024
$ws.Cells.Item(31, 13).Value = -775.5999999999999
$ws.Cells.Item(31, 14).Value = -2337.3024
# Row 34
$ws.Cells.Item(34, 8).Value = 1443.6538
$ws.Cells.Item(34, 9).Value = 1070.6
$ws.Cells.Item(34, 10).Value = 1747.3024
$ws.Cells.Item(34, 11).Value = 1070.6
$ws.Cells.Item(34, 12).Value = 1747.3024
$ws.Cells.Item(34, 13).Value = -868.5999999999999
$ws.Cells.Item(34, 14).Value = -2151.3024
# Row 96
$ws.Cells.Item(96, 8).Value = 18117.75
$ws.Cells.Item(96, 10).Value = 18117.75
$ws.Cells.Item(96, 12).Value = 18117.75
$ws.Cells.Item(96, 14).Value = -23609.75
# Row 109
$ws.Cells.Item(109, 8).Value = 39666.332
$ws.Cells.Item(109, 10).Value = 39666.332
$ws.Cells.Item(109, 12).Value = 39666.332
$ws.Cells.Item(109, 14).Value = -41746.332
# Row 132
$ws.Cells.Item(132, 8).Value = 2686.6667
$ws.Cells.Item(132, 9).Value = 1937.8695
$ws.Cells.Item(132, 10).Value = 5147
$ws.Cells.Item(132, 11).Value = 5813.6085
$ws.Cells.Item(132, 12).Value = 15441
$ws.Cells.Item(132, 13).Value = -3283.6085
$ws.Cells.Item(132, 14).Value = -20501
# Row 134
$ws.Cells.Item(134, 8).Value = 3267.44
$ws.Cells.Item(134, 9).Value = 1359.0769
$ws.Cells.Item(134, 10).Value = 5334.8335
$ws.Cells.Item(134, 11).Value = 4077.2307
$ws.Cells.Item(134, 12).Value = 16004.5005
$ws.Cells.Item(134, 13).Value = -1542.2307
$ws.Cells.Item(134, 14).Value = -21074.5005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 62
$ws.Cells.Item(62, 8).Value = 1000
$ws.Cells.Item(62, 10).Value = 1000
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 14).Value = -4372
# Row 65
$ws.Cells.Item(65, 8).Value = 1000
$ws.Cells.Item(65, 10).Value = 1000
$ws.Cells.Item(65, 12).Value = 9000
$ws.Cells.Item(65, 14).Value = -15864
# Row 68
$ws.Cells.Item(68, 8).Value = 910.131
$ws.Cells.Item(68, 9).Value = 801.78845
$ws.Cells.Item(68, 10).Value = 1086.1875
$ws.Cells.Item(68, 11).Value = 2405.36535
$ws.Cells.Item(68, 12).Value = 3258.5625
$ws.Cells.Item(68, 13).Value = -1594.36535
$ws.Cells.Item(68, 14).Value = -4880.5625
# Row 71
$ws.Cells.Item(71, 8).Value = 910.131
$ws.Cells.Item(71, 9).Value = 801.78845
$ws.Cells.Item(71, 10).Value = 1086.1875
$ws.Cells.Item(71, 11).Value = 7216.09605
$ws.Cells.Item(71, 12).Value = 9775.6875
$ws.Cells.Item(71, 13).Value = -3160.09605
$ws.Cells.Item(71, 14).Value = -17887.6875
# Row 80
$ws.Cells.Item(80, 8).Value = 1011.1111
$ws.Cells.Item(80, 10).Value = 1114
$ws.Cells.Item(80, 12).Value = 3342
$ws.Cells.Item(80, 14).Value = -5214
# Row 83
$ws.Cells.Item(83, 8).Value = 1011.1111
$ws.Cells.Item(83, 10).Value = 1114
$ws.Cells.Item(83, 12).Value = 10026
$ws.Cells.Item(83, 14).Value = -19386
# Row 86
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
# Row 87
$ws.Cells.Item(87, 8).Value = 6260
$ws.Cells.Item(87, 9).Value = 3766.6667
$ws.Cells.Item(87, 10).Value = 10000
$ws.Cells.Item(87, 11).Value = 11300.0001
$ws.Cells.Item(87, 12).Value = 30000
$ws.Cells.Item(87, 13).Value = -10052.0001
$ws.Cells.Item(87, 14).Value = -32496
# Row 89
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
# Row 90
$ws.Cells.Item(90, 8).Value = 6260
$ws.Cells.Item(90, 9).Value = 3766.6667
$ws.Cells.Item(90, 10).Value = 10000
$ws.Cells.Item(90, 11).Value = 33900.0003
$ws.Cells.Item(90, 12).Value = 90000
$ws.Cells.Item(90, 13).Value = -27660.0003
$ws.Cells.Item(90, 14).Value = -102480
# Row 92
$ws.Cells.Item(92, 8).Value = 816.625
$ws.Cells.Item(92, 9).Value = 600
$ws.Cells.Item(92, 10).Value = 847.5714
$ws.Cells.Item(92, 11).Value = 1800
$ws.Cells.Item(92, 12).Value = 2542.7142
$ws.Cells.Item(92, 13).Value = -552
$ws.Cells.Item(92, 14).Value = -5038.7142
# Row 98
$ws.Cells.Item(98, 8).Value = 600
$ws.Cells.Item(98, 9).Value = 466.66666
$ws.Cells.Item(98, 10).Value = 1000
$ws.Cells.Item(98, 11).Value = 1399.99998
$ws.Cells.Item(98, 12).Value = 3000
$ws.Cells.Item(98, 13).Value = 98.00001999999995
$ws.Cells.Item(98, 14).Value = -5996
# Row 99
$ws.Cells.Item(99, 8).Value = 2611.111
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 2611.111
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 7833.333
$ws.Cells.Item(99, 14).Value = -12325.333
$ws.Cells.Item(99, 13).ClearContents()
# Row 101
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 14).ClearContents()
# Row 102
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
# Row 104
$ws.Cells.Item(104, 8).Value = 10000
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 10000
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 30000
$ws.Cells.Item(104, 14).Value = -35242
$ws.Cells.Item(104, 13).ClearContents()
# Row 105
$ws.Cells.Item(105, 8).Value = 10000
$ws.Cells.Item(105, 10).Value = 10000
$ws.Cells.Item(105, 12).Value = 30000
$ws.Cells.Item(105, 14).Value = -35242
# Row 107
$ws.Cells.Item(107, 8).Value = 998.8387
$ws.Cells.Item(107, 9).Value = 325.33334
$ws.Cells.Item(107, 10).Value = 1343.8049
$ws.Cells.Item(107, 11).Value = 976.0000200000001
$ws.Cells.Item(107, 12).Value = 4031.4147
$ws.Cells.Item(107, 13).Value = 943.9999799999999
$ws.Cells.Item(107, 14).Value = -7871.4147

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Cells.Item(40, 8).Value = 3498.6667
$ws.Cells.Item(40, 9).Value = 1811.4
$ws.Cells.Item(40, 11).Value = 1811.4
$ws.Cells.Item(40, 13).Value = -1675.4
# Row 94
$ws.Cells.Item(94, 8).Value = 22165
$ws.Cells.Item(94, 10).Value = 22165
$ws.Cells.Item(94, 12).Value = 22165
$ws.Cells.Item(94, 14).Value = -23517
